# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
# Change cell B11 on the active sheet from "R40" to the text value "1".
# The leading apostrophe forces Excel to store "1" as literal text
# (a shared string) instead of auto-converting it to the number 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "'1"
